$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.248.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.428.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.25%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'414.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.09%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'129.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.53%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.03%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'42.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'9.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.966.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000214"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.00%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'20.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.39%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.439.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'62.207.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.37%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'467.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.23%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'91.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.40%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'13.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.13%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'10.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +20.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'33.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'4.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'11.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.18%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.46%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.167"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.42%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.81%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'40.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.88%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.05%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'58.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +10.36%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0488"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +4.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.325"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'LidoDAOToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'3.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Stellar"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.134"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'145.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.73%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +10.40%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.22%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +19.86%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'16.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.46%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'22.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.79%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₃0516"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +24.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'110.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.83%  "
$ws.Range("E51").Style = "Normal"
